$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; everything currently at/after row 28
# shifts down by one (old row 28 -> new row 29, ..., old row 53 -> new row 54).
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44579
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = 24000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 24500
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Provincia de Diguillín"
$ws.Range("P28").Value = 980
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
